$d = $word.ActiveDocument

$replacements = @(
    @("61×53=", "27×65="),
    @("83×50=", "98×11="),
    @("84×43=", "37×47="),
    @("94×73=", "88×21="),
    @("71×72=", "91×27="),
    @("68×49=", "40×30="),
    @("60×81=", "79×63="),
    @("62×22=", "84×97="),
    @("65×86=", "57×28="),
    @("94×43=", "94×87="),
    @("92×68=", "65×18="),
    @("54×32=", "20×72="),
    @("27×87=", "73×88="),
    @("39×44=", "59×93="),
    @("96×66=", "99×25="),
    @("55×22=", "78×47="),
    @("16×33=", "30×98="),
    @("75×31=", "41×76="),
    @("74×89=", "33×62="),
    @("83×69=", "32×90="),
    @("56×34=", "65×41="),
    @("23×64=", "89×86="),
    @("79×95=", "71×25="),
    @("45×17=", "65×94="),
    @("48×19=", "94×85=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
